$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 996.6667
$ws.Range("I2").Value = 996.6667
$ws.Range("K2").Value = 996.6667
$ws.Range("M2").Value = -883.6667
$ws.Range("H29").Value = 4393.8
$ws.Range("I29").Value = 109
$ws.Range("K29").Value = 327
$ws.Range("M29").Value = -46
$ws.Range("H40").Value = 7555.5557
$ws.Range("I40").Value = 3500
$ws.Range("K40").Value = 3500
$ws.Range("M40").Value = -3325
$ws.Range("H41").Value = 3944.2727
$ws.Range("J41").Value = 4688.2
$ws.Range("L41").Value = 4688.2
$ws.Range("N41").Value = -5568.2
$ws.Range("H76").Value = 2898.8
$ws.Range("I76").Value = 2898.8
$ws.Range("K76").Value = 2898.8
$ws.Range("M76").Value = -2583.8
$ws.Range("H79").Value = 2898.8
$ws.Range("I79").Value = 2898.8
$ws.Range("K79").Value = 2898.8
$ws.Range("M79").Value = -1806.8
$ws.Range("H92").Value = 407.77777
$ws.Range("I92").Value = 407.77777
$ws.Range("K92").Value = 407.77777
$ws.Range("M92").Value = 840.2222300000001
$ws.Range("H95").Value = 52601.25
$ws.Range("J95").Value = 52601.25
$ws.Range("L95").Value = 52601.25
$ws.Range("N95").Value = -58093.25
$ws.Range("H96").Value = 3208
$ws.Range("I96").Value = 2987.7144
$ws.Range("J96").Value = 4750
$ws.Range("K96").Value = 8963.143199999999
$ws.Range("L96").Value = 14250
$ws.Range("M96").Value = -7590.143199999999
$ws.Range("N96").Value = -16996
$ws.Range("H100").Value = 5969.5
$ws.Range("I100").Value = 4997.5
$ws.Range("J100").Value = 6941.5
$ws.Range("K100").Value = 4997.5
$ws.Range("L100").Value = 6941.5
$ws.Range("M100").Value = -4456.5
$ws.Range("N100").Value = -8023.5
$ws.Range("H108").Value = 99995
$ws.Range("J108").Value = 99995
$ws.Range("L108").Value = 99995
$ws.Range("N108").Value = -107675
$ws.Range("H132").Value = 2036.4348
$ws.Range("I132").Value = 2036.4348
$ws.Range("K132").Value = 6109.3044
$ws.Range("M132").Value = -3579.3044
$ws.Range("H137").Value = 2100
$ws.Range("I137").Value = 2133.3333
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 6399.999899999999
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -3849.999899999999
$ws.Range("N137").Value = -11100

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2333
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 2499.5
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 2499.5
$ws.Range("M61").Value = -1788
$ws.Range("N61").Value = -2923.5
$ws.Range("H110").Value = 617
$ws.Range("I110").Value = 617
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 617
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1428
$ws.Range("N110").ClearContents()
$ws.Range("H136").Value = 2333
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 2499.5
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 7498.5
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -12598.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 26838.375
$ws.Range("I86").Value = 2425
$ws.Range("J86").Value = 51251.75
$ws.Range("K86").Value = 2425
$ws.Range("L86").Value = 51251.75
$ws.Range("M86").Value = -1302
$ws.Range("N86").Value = -53497.75
$ws.Range("H89").Value = 26838.375
$ws.Range("I89").Value = 2425
$ws.Range("J89").Value = 51251.75
$ws.Range("K89").Value = 12125
$ws.Range("L89").Value = 256258.75
$ws.Range("M89").Value = -6509
$ws.Range("N89").Value = -267490.75
$ws.Range("H94").Value = 4625
$ws.Range("J94").Value = 6500
$ws.Range("L94").Value = 6500
$ws.Range("N94").Value = -7402
$ws.Range("H99").Value = 2499.8572
$ws.Range("I99").Value = 2249.8333
$ws.Range("K99").Value = 2249.8333
$ws.Range("M99").Value = -751.8332999999998
$ws.Range("H105").Value = 6498.5
$ws.Range("I105").Value = 6498.5
$ws.Range("K105").Value = 6498.5
$ws.Range("M105").Value = -4751.5
$ws.Range("H124").Value = 99995
$ws.Range("J124").Value = 99995
$ws.Range("L124").Value = 99995
$ws.Range("N124").Value = -109815
$ws.Range("H134").Value = 6024.75
$ws.Range("I134").Value = 5799.5
$ws.Range("J134").Value = 6250
$ws.Range("K134").Value = 17398.5
$ws.Range("L134").Value = 18750
$ws.Range("M134").Value = -14863.5
$ws.Range("N134").Value = -23820

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1033.3334
$ws.Range("I23").Value = 950
$ws.Range("J23").Value = 1200
$ws.Range("K23").Value = 950
$ws.Range("L23").Value = 1200
$ws.Range("M23").Value = -710
$ws.Range("N23").Value = -1680
$ws.Range("H27").Value = 1033.3334
$ws.Range("I27").Value = 950
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 950
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = -758
$ws.Range("N27").Value = -1584
$ws.Range("H41").Value = 6000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H99").Value = 800
$ws.Range("J99").Value = 800
$ws.Range("L99").Value = 800
$ws.Range("N99").Value = -3796
$ws.Range("H126").Value = 800
$ws.Range("J126").Value = 800
$ws.Range("L126").Value = 2400
$ws.Range("N126").Value = -7340
$ws.Range("H132").Value = 4800.6
$ws.Range("I132").Value = 4800.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14401.8
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11871.8
$ws.Range("N132").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2425
$ws.Range("I130").Value = 2425
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 7275
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -2255
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 1473.75
$ws.Range("J132").Value = 1648.3334
$ws.Range("L132").Value = 14835.0006
$ws.Range("N132").Value = -19895.0006
$ws.Range("H137").Value = 3197.4
$ws.Range("I137").Value = 1999.5
$ws.Range("K137").Value = 5998.5
$ws.Range("M137").Value = -898.5
$ws.Range("H138").Value = 2974.4
$ws.Range("I138").Value = 3082.6667
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 9248.000100000001
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = -4108.000100000001
$ws.Range("N138").Value = -16280

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 7000
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -31900

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2965.5
$ws.Range("I22").Value = 2360
$ws.Range("J22").Value = 3167.3333
$ws.Range("K22").Value = 2360
$ws.Range("L22").Value = 3167.3333
$ws.Range("M22").Value = -2065
$ws.Range("N22").Value = -3757.3333
$ws.Range("H27").Value = 2965.5
$ws.Range("I27").Value = 2360
$ws.Range("J27").Value = 3167.3333
$ws.Range("K27").Value = 2360
$ws.Range("L27").Value = 3167.3333
$ws.Range("M27").Value = -2253
$ws.Range("N27").Value = -3381.3333
$ws.Range("H32").Value = 4250
$ws.Range("I32").Value = 4250
$ws.Range("K32").Value = 4250
$ws.Range("M32").Value = -3933
$ws.Range("H68").Value = 8135.5557
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 16406.666
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 16406.666
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -17904.666
$ws.Range("H71").Value = 8135.5557
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 16406.666
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 82033.33
$ws.Range("M71").Value = -16256
$ws.Range("N71").Value = -89521.33
$ws.Range("H74").Value = 90000
$ws.Range("I74").Value = 90000
$ws.Range("K74").Value = 90000
$ws.Range("M74").Value = -89002
$ws.Range("H77").Value = 90000
$ws.Range("I77").Value = 90000
$ws.Range("K77").Value = 270000
$ws.Range("M77").Value = -265008
$ws.Range("H100").Value = 3355.111
$ws.Range("I100").Value = 3149.25
$ws.Range("J100").Value = 3519.8
$ws.Range("K100").Value = 3149.25
$ws.Range("L100").Value = 3519.8
$ws.Range("M100").Value = -2608.25
$ws.Range("N100").Value = -4601.8
$ws.Range("H127").Value = 90354.5
$ws.Range("J127").Value = 90354.5
$ws.Range("L127").Value = 90354.5
$ws.Range("N127").Value = -100274.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H70").Value = 90000
$ws.Range("I70").Value = 90000
$ws.Range("J70").Value = 90000
$ws.Range("K70").Value = 90000
$ws.Range("L70").Value = 90000
$ws.Range("M70").Value = -89685
$ws.Range("N70").Value = -90630
$ws.Range("H73").Value = 90000
$ws.Range("I73").Value = 90000
$ws.Range("J73").Value = 90000
$ws.Range("K73").Value = 90000
$ws.Range("L73").Value = 90000
$ws.Range("M73").Value = -88908
$ws.Range("N73").Value = -92184
$ws.Range("H107").Value = 1929.8334
$ws.Range("I107").Value = 915.8
$ws.Range("K107").Value = 2747.4
$ws.Range("M107").Value = -827.3999999999996
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H136").Value = 7092.6
$ws.Range("I136").Value = 8000
$ws.Range("J136").Value = 6865.75
$ws.Range("K136").Value = 24000
$ws.Range("L136").Value = 20597.25
$ws.Range("M136").Value = -21450
$ws.Range("N136").Value = -25697.25
